$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.06428691266220485
$ws.Range("E2").Value = 0.06428691266220485

# Row 3
$ws.Range("D3").Value = 0.0004177273670890082
$ws.Range("E3").Value = 0.0004177273670890082

# Row 4
$ws.Range("D4").Value = 0.9790669880789175
$ws.Range("E4").Value = 0.9790669880789175

# Row 5
$ws.Range("D5").Value = 0.000399177427840021
$ws.Range("E5").Value = 0.000399177427840021

# Row 6
$ws.Range("D6").Value = 0.05219987094211664
$ws.Range("E6").Value = 0.05219987094211664

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"2.419878922030138E-09"
$ws.Range("E8").Value = 0.999999997580121

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.009098934946045283
$ws.Range("E9").Value = 0.9909010650539547

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"3.767160850072398E-06"
$ws.Range("E10").Value = 0.99999623283915

# Row 11
$ws.Range("D11").Value = 0.9999999950982748
$ws.Range("E11").Value = [double]"4.901725247918876E-09"
$ws.Range("F11").Value = 4.101563453674316
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.02979376804258702
$ws.Range("E12").Value = 0.02979376804258702

# Row 13
$ws.Range("D13").Value = [double]"8.606573005312807E-05"
$ws.Range("E13").Value = [double]"8.606573005312807E-05"

# Row 14
$ws.Range("D14").Value = 0.9961459094070299
$ws.Range("E14").Value = 0.9961459094070299

# Row 15
$ws.Range("D15").Value = 0.001133821599325619
$ws.Range("E15").Value = 0.001133821599325619

# Row 16
$ws.Range("D16").Value = 0.01548560677879582
$ws.Range("E16").Value = 0.01548560677879582

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"1.463037552608386E-13"
$ws.Range("E18").Value = 0.9999999999998537

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.002845659743289048
$ws.Range("E19").Value = 0.9971543402567109

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"0.0001149059028904656"
$ws.Range("E20").Value = 0.9998850940971096

# Row 21
$ws.Range("D21").Value = 0.9999999999969675
$ws.Range("E21").Value = [double]"3.032463169461153E-12"
$ws.Range("F21").Value = 5.009214401245117
$ws.Range("G21").Value = 0.6
